$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EDA")
$ws.Columns("L").ColumnWidth = 13.33
